$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 4155.125
$ws.Range("I53").Value = 316.66666
$ws.Range("J53").Value = 6458.2
$ws.Range("K53").Value = 316.66666
$ws.Range("L53").Value = 6458.2
$ws.Range("M53").Value = 320.33334
$ws.Range("N53").Value = -7732.2
$ws.Range("H92").Value = 533.2692
$ws.Range("I92").Value = 518.15
$ws.Range("J92").Value = 583.6667
$ws.Range("K92").Value = 518.15
$ws.Range("L92").Value = 583.6667
$ws.Range("M92").Value = 729.85
$ws.Range("N92").Value = -3079.6667
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H98").Value = 383.42856
$ws.Range("I98").Value = 384
$ws.Range("J98").Value = 380
$ws.Range("K98").Value = 384
$ws.Range("L98").Value = 380
$ws.Range("M98").Value = 1114
$ws.Range("N98").Value = -3376
$ws.Range("H100").Value = 2522.2222
$ws.Range("J100").Value = 3350
$ws.Range("L100").Value = 3350
$ws.Range("N100").Value = -4432
$ws.Range("H112").Value = 3907295.2
$ws.Range("J112").Value = 1098.2333
$ws.Range("L112").Value = 3294.699900000001
$ws.Range("N112").Value = -5510.699900000001
$ws.Range("H121").Value = 892.12
$ws.Range("J121").Value = 892.12
$ws.Range("L121").Value = 2676.36
$ws.Range("N121").Value = -6170.360000000001
$ws.Range("H122").Value = 383.42856
$ws.Range("I122").Value = 384
$ws.Range("J122").Value = 380
$ws.Range("K122").Value = 1152
$ws.Range("L122").Value = 1140
$ws.Range("M122").Value = 1298
$ws.Range("N122").Value = -6040
$ws.Range("H135").Value = 27786502
$ws.Range("I135").Value = 1043.4
$ws.Range("J135").Value = 62518330
$ws.Range("K135").Value = 9390.6
$ws.Range("L135").Value = 562664970
$ws.Range("M135").Value = -6855.6
$ws.Range("N135").Value = -562670040
$ws.Range("H138").Value = 2549.077
$ws.Range("J138").Value = 3043.4517
$ws.Range("L138").Value = 9130.355100000001
$ws.Range("N138").Value = -19410.3551
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23022.94
$ws.Range("I32").Value = 26918.324
$ws.Range("J32").Value = 7441.4
$ws.Range("K32").Value = 26918.324
$ws.Range("L32").Value = 7441.4
$ws.Range("M32").Value = -26631.324
$ws.Range("N32").Value = -8015.4
$ws.Range("H101").Value = 40250
$ws.Range("J101").Value = 40250
$ws.Range("L101").Value = 40250
$ws.Range("N101").Value = -46740
$ws.Range("H110").Value = 2560.4546
$ws.Range("I110").Value = 2515
$ws.Range("J110").Value = 2615
$ws.Range("K110").Value = 2515
$ws.Range("L110").Value = 2615
$ws.Range("M110").Value = -470
$ws.Range("N110").Value = -6705
$ws.Range("H122").Value = 2188.0454
$ws.Range("J122").Value = 2416.5
$ws.Range("L122").Value = 7249.5
$ws.Range("N122").Value = -12149.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1128.1666
$ws.Range("I64").Value = 3166.6667
$ws.Range("J64").Value = 448.66666
$ws.Range("K64").Value = 3166.6667
$ws.Range("L64").Value = 448.66666
$ws.Range("M64").Value = -2941.6667
$ws.Range("N64").Value = -898.66666
$ws.Range("H67").Value = 1128.1666
$ws.Range("I67").Value = 3166.6667
$ws.Range("J67").Value = 448.66666
$ws.Range("K67").Value = 3166.6667
$ws.Range("L67").Value = 448.66666
$ws.Range("M67").Value = -2386.6667
$ws.Range("N67").Value = -2008.66666
$ws.Range("H80").Value = 772.04346
$ws.Range("J80").Value = 477.70587
$ws.Range("L80").Value = 477.70587
$ws.Range("N80").Value = -2473.70587
$ws.Range("H83").Value = 772.04346
$ws.Range("J83").Value = 477.70587
$ws.Range("L83").Value = 2388.52935
$ws.Range("N83").Value = -12372.52935
$ws.Range("H86").Value = 27873.947
$ws.Range("I86").Value = 39662
$ws.Range("J86").Value = 2333.1667
$ws.Range("K86").Value = 39662
$ws.Range("L86").Value = 2333.1667
$ws.Range("M86").Value = -38539
$ws.Range("N86").Value = -4579.1667
$ws.Range("H89").Value = 27873.947
$ws.Range("I89").Value = 39662
$ws.Range("J89").Value = 2333.1667
$ws.Range("K89").Value = 198310
$ws.Range("L89").Value = 11665.8335
$ws.Range("M89").Value = -192694
$ws.Range("N89").Value = -22897.8335
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 133.5
$ws.Range("J7").Value = 240
$ws.Range("L7").Value = 240
$ws.Range("N7").Value = -466
$ws.Range("H132").Value = 69302.125
$ws.Range("I132").Value = 500001
$ws.Range("J132").Value = 7773.7144
$ws.Range("K132").Value = 1500003
$ws.Range("L132").Value = 23321.1432
$ws.Range("M132").Value = -1497473
$ws.Range("N132").Value = -28381.1432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2939.5
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 3252.6667
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 9758.000100000001
$ws.Range("M64").Value = -5730
$ws.Range("N64").Value = -10298.0001
$ws.Range("H67").Value = 2939.5
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 3252.6667
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 9758.000100000001
$ws.Range("M67").Value = -5064
$ws.Range("N67").Value = -11630.0001
$ws.Range("H68").Value = 1227.1923
$ws.Range("I68").Value = 433.66666
$ws.Range("J68").Value = 1330.6957
$ws.Range("K68").Value = 1300.99998
$ws.Range("L68").Value = 3992.0871
$ws.Range("M68").Value = -489.9999800000001
$ws.Range("N68").Value = -5614.0871
$ws.Range("H71").Value = 1227.1923
$ws.Range("I71").Value = 433.66666
$ws.Range("J71").Value = 1330.6957
$ws.Range("K71").Value = 3902.99994
$ws.Range("L71").Value = 11976.2613
$ws.Range("M71").Value = 153.0000600000003
$ws.Range("N71").Value = -20088.2613
$ws.Range("H107").Value = 3703.3438
$ws.Range("I107").Value = 7629.643
$ws.Range("J107").Value = 649.55554
$ws.Range("K107").Value = 22888.929
$ws.Range("L107").Value = 1948.66662
$ws.Range("M107").Value = -20968.929
$ws.Range("N107").Value = -5788.66662
$ws.Range("H129").Value = 278401.22
$ws.Range("I129").Value = 469.1111
$ws.Range("J129").Value = 556333.3
$ws.Range("K129").Value = 1407.3333
$ws.Range("L129").Value = 1668999.9
$ws.Range("M129").Value = 3592.6667
$ws.Range("N129").Value = -1678999.9
$ws.Range("H131").Value = 127425.836
$ws.Range("J131").Value = 135969.56
$ws.Range("L131").Value = 407908.68
$ws.Range("N131").Value = -417988.68
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10732.667
$ws.Range("I80").Value = 17272.857
$ws.Range("J80").Value = 5010
$ws.Range("K80").Value = 17272.857
$ws.Range("L80").Value = 5010
$ws.Range("M80").Value = -16274.857
$ws.Range("N80").Value = -7006
$ws.Range("H83").Value = 10732.667
$ws.Range("I83").Value = 17272.857
$ws.Range("J83").Value = 5010
$ws.Range("K83").Value = 86364.285
$ws.Range("L83").Value = 25050
$ws.Range("M83").Value = -81372.285
$ws.Range("N83").Value = -35034
$ws.Range("H113").Value = 3500
$ws.Range("I113").Value = 2931.25
$ws.Range("J113").Value = 6533.3335
$ws.Range("K113").Value = 2931.25
$ws.Range("L113").Value = 6533.3335
$ws.Range("M113").Value = -761.25
$ws.Range("N113").Value = -10873.3335
$ws.Range("H122").Value = 3541.0952
$ws.Range("I122").Value = 2558.1538
$ws.Range("J122").Value = 5138.375
$ws.Range("K122").Value = 7674.4614
$ws.Range("L122").Value = 15415.125
$ws.Range("M122").Value = -5224.4614
$ws.Range("N122").Value = -20315.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3493.2222
$ws.Range("I7").Value = 4190.625
$ws.Range("J7").Value = 2478.818
$ws.Range("K7").Value = 4190.625
$ws.Range("L7").Value = 2478.818
$ws.Range("M7").Value = -4078.625
$ws.Range("N7").Value = -2702.818
$ws.Range("H55").Value = 252.82608
$ws.Range("I55").Value = 186.45454
$ws.Range("J55").Value = 313.66666
$ws.Range("K55").Value = 186.45454
$ws.Range("L55").Value = 313.66666
$ws.Range("M55").Value = -13.45454000000001
$ws.Range("N55").Value = -659.66666
$ws.Range("H122").Value = 2504.16
$ws.Range("I122").Value = 2067
$ws.Range("K122").Value = 6201
$ws.Range("M122").Value = -3751
$ws.Range("H126").Value = 3493.2222
$ws.Range("I126").Value = 4190.625
$ws.Range("J126").Value = 2478.818
$ws.Range("K126").Value = 12571.875
$ws.Range("L126").Value = 7436.454000000001
$ws.Range("M126").Value = -10101.875
$ws.Range("N126").Value = -12376.454
$ws.Range("H132").Value = 863085.4
$ws.Range("I132").Value = 1722927.4
$ws.Range("J132").Value = 3243.2856
$ws.Range("K132").Value = 5168782.199999999
$ws.Range("L132").Value = 9729.856800000001
$ws.Range("M132").Value = -5166252.199999999
$ws.Range("N132").Value = -14789.8568
$ws.Range("H136").Value = 26096.857
$ws.Range("I136").Value = 32502.125
$ws.Range("K136").Value = 97506.375
$ws.Range("M136").Value = -94956.375
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1438.1666
$ws.Range("I132").Value = 1138.12
$ws.Range("K132").Value = 3414.36
$ws.Range("M132").Value = -884.3599999999997
